# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.046.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.925.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.75%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5163"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4002"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08469"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.122"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.320"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.921.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.359"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001114"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06752"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.050"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.058.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.208"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.143.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.461"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.077"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1057"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.081"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.665"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02498"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06599"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2212"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.242"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.018"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.188"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6527"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.238"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6134"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.722"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.054"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.58%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.00%  "
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.239"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.150"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.82%  "
